# Adds a new "p2p BitTorrent" protocol column-pair (F:G) to each of the three
# worksheets, pushing the former "HTTP/1.1" column-pair (B:C) data into the
# new F:G slot and replacing B:C with fresh p2p BitTorrent measurements for
# the 1MB / 10MB rows (the 10kB / 100kB rows have no p2p data, so they go
# blank).
#
# xlPasteFormats = -4122 (used instead of a full Range.Copy into the merged
# header cells, which would otherwise trigger Excel's "outer border on
# merge" auto-formatting and mint new style records).

$wb = $excel.ActiveWorkbook

# New row 6 / row 7 values for column B/C (p2p BitTorrent), per sheet index.
$newB6 = @(2.349340756734212, 10114.69803173217, 1.00267129492618)
$newC6 = @(0.6341095253989198, 3435.038111728009, 0.00014128277808102)
$newB7 = @(4.431923389434814, 56635.20281022038, 1.00177945368595)

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- 1. Make room: give F:G the same column width as the other cols ---
    # (14.17 "characters" is what Excel's ColumnWidth model round-trips to
    # the stored width="15" units used by columns A:E in this sheet)
    $ws.Columns.Item(6).ColumnWidth = 14.17
    $ws.Columns.Item(7).ColumnWidth = 14.17

    # --- 2. Row 1 header: move the merged "HTTP/1.1" header into F1:G1 ---
    # Flip on the merge first (cheap on still-unstyled cells, no border
    # side-effects), THEN paste formats into each half separately so both
    # halves keep plain style 1 (matching B1/C1's un-merged-yet look).
    $ws.Range("F1:G1").MergeCells = $true
    $ws.Range("B1").Copy()
    $ws.Range("F1").PasteSpecial(-4122)
    $ws.Range("C1").Copy()
    $ws.Range("G1").PasteSpecial(-4122)
    $ws.Range("F1").Value = "HTTP/1.1"

    # --- 3. Row 2 (Mean / Std Dev labels) ---
    $ws.Range("B2:C2").Copy($ws.Range("F2:G2"))

    # --- 4. Row 3 (blank bordered spacer row) ---
    $ws.Range("B3:C3").Copy($ws.Range("F3:G3"))

    # --- 5. Rows 4:5 (10kB / 100kB data) ---
    $ws.Range("B4:C5").Copy($ws.Range("F4:G5"))

    # --- 6. Rows 6:7 (1MB / 10MB totals, unstyled cells) ---
    $ws.Range("B6:C7").Copy($ws.Range("F6:G7"))

    # --- 7. Relabel the original header: HTTP/1.1 -> p2p BitTorrent ---
    $ws.Range("B1").Value = "p2p BitTorrent"

    # --- 8. p2p BitTorrent has no 10kB/100kB samples: blank out B4:C5 ---
    $ws.Range("B3").Copy()
    $ws.Range("B4:C5").PasteSpecial(-4122)
    $ws.Range("B4:C5").ClearContents()

    # --- 9. Write the new p2p BitTorrent totals for 1MB / 10MB ---
    $idx = $i - 1
    $ws.Range("B6").Value = $newB6[$idx]
    $ws.Range("C6").Value = $newC6[$idx]
    $ws.Range("B7").Value = $newB7[$idx]
    # C7 is unchanged (stays 0)
}
